$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the D2 diode row: change Value (B3) and P/N (D3)
$ws.Range("B3").Value = "M7"
$ws.Range("D3").Value = "RS1MWF-7"

# Update the active selection to match the saved view state
$ws.Range("I11").Select()
